$wb = $excel.ActiveWorkbook

# --- Sheet "Table-3.1" (3rd sheet): fill in top/second/third English-speaking country answers ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C5").Value = "USA"
$ws3.Range("C6").Value = "Other Countries"
$ws3.Range("C7").Value = "GBR"

# --- Sheet "Table-5.1" (4th sheet): fill in sector-wise investment answers ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C5").Value = 12150
$ws4.Range("C6").Value = 108531347515
$ws4.Range("C7").Value = "Others"
$ws4.Range("C8").Value = "Cleantech & Semiconductors"
$ws4.Range("C9").Value = "Social & Finance & Analytics & Advertising"
$ws4.Range("C10").Value = 2923
$ws4.Range("C11").Value = 2297
$ws4.Range("C12").Value = 1912

# Column C grew much wider text, column D is effectively empty now -- resize
# (best-fit) both columns to match the new content.
$ws4.Columns.Item(3).AutoFit()
$ws4.Columns.Item(4).AutoFit()

# Rows 4 and 6 go back to the default (auto) row height, row 5 keeps an
# explicit 15pt height.
$ws4.Rows.Item(4).AutoFit()
$ws4.Rows.Item(5).RowHeight = 15
$ws4.Rows.Item(6).AutoFit()

# Select a cell on Table-3.1 (its new active cell) ...
$ws3.Range("C7").Select() | Out-Null

# ... then make Table-5.1 the active sheet with C14 selected, so it becomes
# the workbook's active tab (and Table - 2.1 loses its previous tabSelected).
$ws4.Range("C14").Select() | Out-Null
